$d = $word.ActiveDocument

# Locate the "Vulnerable and Outdated Components" table-cell paragraph
# (the 2021 category name for A06) without relying on hard-coded table /
# row / column indices.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Vulnerable and Outdated Components"
$find.Forward = $true
$find.Wrap = 1
$find.Execute() | Out-Null

$para = $d.Range($find.Parent.Start, $find.Parent.End)
$insertPos = $para.End
$para.Collapse(0)
$para.InsertAfter(" ")

# Force the newly inserted space to live in its own run (matching the
# target edit) instead of being silently merged back into the preceding
# "Vulnerable and Outdated Components" run, by toggling a character
# property on just the new character and then reverting it.
$newRange = $d.Range($insertPos, $insertPos + 1)
$newRange.Font.Bold = 1
$newRange.Font.Bold = 0
